# 2.4.1 Pause works in --headed mode.
#
# This script:
#  1. Makes the "attrib_href" sheet visible again (was hidden).
#  2. Adjusts selections on a couple of sheets (attrib_href, dblclick) that
#     were touched as part of the commit's re-save.
#  3. Appends a brand-new "pause" worksheet (sheetId 11) after "file",
#     modeled on the existing "dblclick" action-sheet, describing the new
#     pause test case, and makes it the active sheet/tab (matching
#     activeTab="10" in the workbook view).

$wb = $excel.ActiveWorkbook

# --- 1. Unhide attrib_href -------------------------------------------------
$wsAttrib = $wb.Worksheets.Item("attrib_href")
$wsAttrib.Visible = -1   # xlSheetVisible

# --- 2. Cosmetic selection tweaks on existing sheets -----------------------
# These don't activate/select the sheet itself as "current" for the final
# save - only change the remembered cell selection inside each sheet.
$wsAttrib.Range("C33").Select()

$wsDblclick = $wb.Worksheets.Item("dblclick")
$wsDblclick.Range("A11").Select()

# --- 3. Add the new "pause" sheet at the end -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPause = $wb.Worksheets.Add($null, $lastSheet)
$wsPause.Name = "pause"

# Header row (same header used by every case-sheet in this workbook).
$wsPause.Cells.Item(1,1).Value = "Desc"
$wsPause.Cells.Item(1,2).Value = "Steps"
$wsPause.Cells.Item(1,3).Value = "Locator"
$wsPause.Cells.Item(1,4).Value = "Action"
$wsPause.Cells.Item(1,5).Value = "Data"

# Case description row.
$wsPause.Cells.Item(2,1).Value = "action – pause and print"

# Step 1: navigate + assert title (mirrors the other action sheets).
$wsPause.Cells.Item(3,2).Value = "Case 1"
$wsPause.Cells.Item(3,3).Value = "https://tecagile.com/double-click-test/"
$wsPause.Cells.Item(3,4).Value = "url"

$wsPause.Cells.Item(4,4).Value = "title"
$wsPause.Cells.Item(4,5).Value = "Online Double"

# Step 2: print a line, pause, then print another line.
$wsPause.Cells.Item(5,4).Value = "print"
$wsPause.Cells.Item(5,5).Value = "line before pause"

$wsPause.Cells.Item(6,4).Value = "pause"

$wsPause.Cells.Item(7,4).Value = "print"
$wsPause.Cells.Item(7,5).Value = "line after pause"

# Trailing blank-but-present rows, matching the source sheet's used range
# (dimension A1:E11).
$wsPause.Cells.Item(8,5).Value = ""
$wsPause.Cells.Item(9,5).Value = ""
$wsPause.Cells.Item(10,5).Value = ""
$wsPause.Cells.Item(11,5).Value = ""

# Column widths matching the sibling "dblclick" sheet layout.
$wsPause.Columns.Item(3).ColumnWidth = 32.74
$wsPause.Columns.Item(4).ColumnWidth = 7.61
$wsPause.Columns.Item(5).ColumnWidth = 16.2

# Make "pause" the active sheet/tab (activeTab="10"), with the same
# selection the source workbook recorded for it (F15).
$wsPause.Activate()
$wsPause.Range("F15").Select()

Write-Output "pause sheet added; attrib_href unhidden; selections updated."
